$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.283.21'
$ws.Cells.Item(2, 5).Value = '  +0.16%  '
$ws.Cells.Item(3, 4).Value = '1.841.56'
$ws.Cells.Item(3, 5).Value = '  +0.17%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.9986'
$ws.Cells.Item(4, 5).Value = '  -0.22%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '241.02'
$ws.Cells.Item(5, 5).Value = '  -0.76%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.6693'
$ws.Cells.Item(6, 5).Value = '  -2.27%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.9999'
$ws.Cells.Item(7, 5).Value = '  -0.14%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.07424'
$ws.Cells.Item(8, 5).Value = '  -1.11%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.2963'
$ws.Cells.Item(9, 5).Value = '  -1.94%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '22.92'
$ws.Cells.Item(10, 5).Value = '  -1.08%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07715'
$ws.Cells.Item(11, 5).Value = '  +0.84%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '5.039'
$ws.Cells.Item(12, 5).Value = '  -0.59%  '
$ws.Cells.Item(13, 2).Value = 'Polygon'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.6801'
$ws.Cells.Item(13, 5).Value = '  -0.42%  '
$ws.Cells.Item(14, 2).Value = 'WrappedEther'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(14, 4).Value = '1.773.40'
$ws.Cells.Item(14, 5).Value = '  -3.82%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '86.41'
$ws.Cells.Item(15, 5).Value = '  -3.00%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '6.215'
$ws.Cells.Item(16, 5).Value = '  -1.07%  '
$ws.Cells.Item(17, 4).Value = '29.174.70'
$ws.Cells.Item(17, 5).Value = '  -0.27%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.000008254'
$ws.Cells.Item(18, 5).Value = '  +0.45%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '229.41'
$ws.Cells.Item(19, 5).Value = '  -1.96%  '
$ws.Cells.Item(20, 5).Value = '  -0.07%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '0.9989'
$ws.Cells.Item(21, 5).Value = '  -0.14%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '7.285'
$ws.Cells.Item(22, 5).Value = '  -2.33%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.9996'
$ws.Cells.Item(23, 5).Value = '  -0.20%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '160.29'
$ws.Cells.Item(24, 5).Value = '  +0.37%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '8.733'
$ws.Cells.Item(25, 5).Value = '  -0.88%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.1417'
$ws.Cells.Item(26, 5).Value = '  -2.40%  '
$ws.Cells.Item(27, 5).Value = '  +0.01%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '1.505'
$ws.Cells.Item(28, 5).Value = '  -1.01%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '4.214'
$ws.Cells.Item(29, 5).Value = '  +0.13%  '
$ws.Cells.Item(30, 5).Value = '  -0.71%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.198'
$ws.Cells.Item(31, 5).Value = '  -0.01%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.05355'
$ws.Cells.Item(32, 5).Value = '  +4.53%  '
$ws.Cells.Item(33, 5).Value = '  -1.17%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.869'
$ws.Cells.Item(34, 5).Value = '  +1.51%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.138'
$ws.Cells.Item(35, 5).Value = '  +0.21%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '2.681'
$ws.Cells.Item(36, 5).Value = '  +0.24%  '
$ws.Cells.Item(37, 4).Value = '1.333.24'
$ws.Cells.Item(37, 5).Value = '  +3.46%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.01803'
$ws.Cells.Item(38, 5).Value = '  -2.24%  '
$ws.Cells.Item(39, 5).Value = '  +1.20%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.9203'
$ws.Cells.Item(40, 5).Value = '  -2.11%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '6.014'
$ws.Cells.Item(41, 5).Value = '  +6.58%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.002'
$ws.Cells.Item(42, 5).Value = '  +0.12%  '
$ws.Cells.Item(43, 5).Value = '  -1.61%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.07957'
$ws.Cells.Item(44, 5).Value = '  +16.02%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.00000000124'
$ws.Cells.Item(45, 5).Value = '  +1.41%  '
$ws.Cells.Item(46, 2).Value = 'Mantle'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.5163'
$ws.Cells.Item(46, 5).Value = '  -0.59%  '
$ws.Cells.Item(47, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(47, 4).Value = '1.950.67'
$ws.Cells.Item(47, 5).Value = '  -2.14%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '63.99'
$ws.Cells.Item(48, 5).Value = '  +2.01%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.767'
$ws.Cells.Item(49, 5).Value = '  +0.31%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '9.261'
$ws.Cells.Item(50, 5).Value = '  -4.06%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.05951'
$ws.Cells.Item(51, 5).Value = '  +0.41%  '
